$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Fecha" (D) and "Volumen" (J) columns to the new weekly values.
$ws.Range("D2").Value = 44651
$ws.Range("D3").Value = 44649
$ws.Range("J3").Value = 60
$ws.Range("D4").Value = 44635
$ws.Range("J4").Value = 100
$ws.Range("D5").Value = 44642
$ws.Range("J5").Value = 100
$ws.Range("D8").Value = 44658
$ws.Range("J8").Value = 80
$ws.Range("D10").Value = 44637
$ws.Range("J10").Value = 100
$ws.Range("D11").Value = 44664
$ws.Range("J11").Value = 160
$ws.Range("D12").Value = 44656
$ws.Range("D13").Value = 44628
$ws.Range("J13").Value = 60
$ws.Range("D14").Value = 44659
$ws.Range("J14").Value = 80
